$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) AUDIENCE paragraph: replace + re-split into many runs + insert bookmark
# ---------------------------------------------------------------------------
$old = "This session is ideal for individuals with basic knowledge in the subject and is appropriate for a beginner in the area."
$new = "The session is ideal for individuals with basic knowledge related to Natural Language Processing (NLP) and Machine Learning. This research is appropriate for any beginner in these areas."

$r = $d.Content
$r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $r.Start
$r.Text = $new

# segment boundaries (relative offsets into $new)
$bounds = @(0,2,3,25,58,102,103,123,139,159,161,176,179,184,185,186)

for ($i = 0; $i -lt $bounds.Length - 1; $i++) {
  $segStart = $paraStart + $bounds[$i]
  $segEnd = $paraStart + $bounds[$i+1]
  $seg = $d.Range($segStart, $segEnd)
  $seg.Bold = 1
  $seg.Bold = 0
}

# bookmark _GoBack goes right after segment index 6 ("and Machine Learning"), i.e. at offset 123
$bmPos = $paraStart + 123
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
